$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.567.22'
$ws.Range('E2').Value = '  -4.19%  '
$ws.Range('D3').Value = '3.345.33'
$ws.Range('E3').Value = '  -2.19%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '572.25'
$ws.Range('E5').Value = '  -0.97%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '132.42'
$ws.Range('E6').Value = '  +4.07%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '3.344.94'
$ws.Range('E8').Value = '  -2.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.476'
$ws.Range('E9').Value = '  -0.38%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.59'
$ws.Range('E10').Value = '  +1.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.121'
$ws.Range('E11').Value = '  -0.85%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.386'
$ws.Range('E12').Value = '  +1.87%  '
$ws.Range('D13').Value = '3.918.35'
$ws.Range('E13').Value = '  -2.18%  '
$ws.Range('E14').Value = '  +0.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000173'
$ws.Range('E15').Value = '  -0.89%  '
$ws.Range('D16').Value = '3.347.67'
$ws.Range('E16').Value = '  -2.13%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '24.94'
$ws.Range('E17').Value = '  +0.25%  '
$ws.Range('D18').Value = '60.719.59'
$ws.Range('E18').Value = '  -3.93%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.88'
$ws.Range('E19').Value = '  +5.17%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.77'
$ws.Range('E20').Value = '  +1.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.22'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '373.66'
$ws.Range('E22').Value = '  -1.48%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.563'
$ws.Range('E23').Value = '  +0.38%  '
$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('B25').Value = 'WrappedeETH'
$ws.Range('C25').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D25').Value = '3.477.74'
$ws.Range('E25').Value = '  -2.30%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '70.06'
$ws.Range('E26').Value = '  -3.90%  '
$ws.Range('E27').Value = '  +5.88%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.68'
$ws.Range('E28').Value = '  +19.79%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.64'
$ws.Range('E29').Value = '  +9.11%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.997'
$ws.Range('E30').Value = '  -0.11%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.07'
$ws.Range('E31').Value = '  +2.02%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.16'
$ws.Range('E32').Value = '  -0.70%  '
$ws.Range('E33').Value = '  +0.62%  '
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('D35').Value = '3.377.02'
$ws.Range('E35').Value = '  -2.20%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '23.09'
$ws.Range('E36').Value = '  +1.17%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.44'
$ws.Range('E37').Value = '  +2.88%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.96'
$ws.Range('E38').Value = '  +3.07%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.54'
$ws.Range('E39').Value = '  +2.25%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '162.51'
$ws.Range('E40').Value = '  -0.99%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0781'
$ws.Range('E41').Value = '  +2.30%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('E43').Value = '  +12.14%  '
$ws.Range('E44').Value = '  +2.72%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '41.40'
$ws.Range('E45').Value = '  -1.21%  '
$ws.Range('E46').Value = '  -4.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '23.44'
$ws.Range('E47').Value = '  +2.13%  '
$ws.Range('E48').Value = '  +0.02%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.93'
$ws.Range('E49').Value = '  +3.32%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '22.83'
$ws.Range('E50').Value = '  +12.07%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.891'
$ws.Range('E51').Value = '  +2.99%  '
